$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Vlan_params")
$ws2 = $wb.Worksheets.Item("Data Types")

# Copy the existing text "True" cell (Data Types!E7) so the new values land
# as text (matching "True"/"False" already on the Data Types sheet) instead
# of being auto-coerced to a boolean by a plain .Value assignment.
$ws2.Range("E7").Copy()
$ws1.Range("B3").PasteSpecial(-4163)  # xlPasteValues

$ws2.Range("E7").Copy()
$ws1.Range("B9").PasteSpecial(-4163)  # xlPasteValues

$excel.CutCopyMode = 0

$ws1.Range("E7").Select()
